# ===== Weather workbook edit script =====
# Applies: captain tags ($ captain suffix) to player cells in RoC/turf/Sheet4,
# fixes Column A (manager name) in "turf" sheet, and appends 3 new manager rows
# (rows 7-9) plus a stray row 10 cell in "RoC", then re-selects RoC as active sheet.

$wb = $excel.ActiveWorkbook

function Set-PlainText($sheet, $cellAddr, [string]$val) {
    # Writes a plain string value (never auto-coerced to a number/date).
    $helper = $sheet.Cells.Item(500, 500)
    $helper.NumberFormat = "@"
    $helper.Value = $val
    $helper.Copy() | Out-Null
    $target = $sheet.Range($cellAddr)
    $target.PasteSpecial(-4163) | Out-Null
    $helper.ClearContents() | Out-Null
    $helper.NumberFormat = "General"
}

# ----- turf sheet -----
$wsTurf = $wb.Worksheets.Item("turf")

# Column A should mirror column F (manager name) - fixes a data-entry bug
Set-PlainText $wsTurf "A2" "Sahil Lampard"
Set-PlainText $wsTurf "A3" "Saurabh Tamang"
Set-PlainText $wsTurf "A4" "Pranesh Sharma"
Set-PlainText $wsTurf "A5" "Sanjeev Kumar"
Set-PlainText $wsTurf "A6" "Kumar Anku"
Set-PlainText $wsTurf "A7" "Kushal Sahota"
Set-PlainText $wsTurf "A8" "Regean Lama"
Set-PlainText $wsTurf "A9" "pranam rai"
Set-PlainText $wsTurf "A10" "Ashim Lama"
Set-PlainText $wsTurf "A11" "Vivek Pradhan"
Set-PlainText $wsTurf "A12" "Sibin _"
Set-PlainText $wsTurf "A13" "Omkar Subba"
Set-PlainText $wsTurf "A14" "Indrajeet Singh"
Set-PlainText $wsTurf "A15" "Rishi Thulung"
Set-PlainText $wsTurf "A16" "Dipen Thapa"
Set-PlainText $wsTurf "A17" "Nishant Singh"
Set-PlainText $wsTurf "A18" "Ashish Subba"
Set-PlainText $wsTurf "A19" "Rohan Rai"
Set-PlainText $wsTurf "A20" "Raj Kanango"
Set-PlainText $wsTurf "A21" "Gaurav Shrestha"
Set-PlainText $wsTurf "A22" "Bishal Rai"
Set-PlainText $wsTurf "A23" "RAJIYUNG Sun"
Set-PlainText $wsTurf "A24" "Ashish Rai"
Set-PlainText $wsTurf "A25" "prabin kumar"
Set-PlainText $wsTurf "A26" "Ravish Verma"
Set-PlainText $wsTurf "A27" "Prayash Thapa"
Set-PlainText $wsTurf "A28" "Tenzing Ninjey"
Set-PlainText $wsTurf "A29" "Harsh Raval"
Set-PlainText $wsTurf "A30" "PRAYAS BAJGAI"
Set-PlainText $wsTurf "A31" "AVIJIT DAS"
Set-PlainText $wsTurf "A32" "Abhijeet Singh"
Set-PlainText $wsTurf "A33" "SATYAM RAI"

# Tag each captain pick with a "$ captain" suffix
Set-PlainText $wsTurf "Z2" "Salah 26`$ captain"
Set-PlainText $wsTurf "AD3" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "AD4" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "X5" "Salah 26`$ captain"
Set-PlainText $wsTurf "AC6" "Antonio 4`$ captain"
Set-PlainText $wsTurf "AB7" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "AD8" "Ronaldo 2`$ captain"
Set-PlainText $wsTurf "AC9" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "AB10" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "AD11" "Ronaldo 2`$ captain"
Set-PlainText $wsTurf "AB12" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "AC13" "Antonio 4`$ captain"
Set-PlainText $wsTurf "AC14" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "AD15" "Ronaldo 2`$ captain"
Set-PlainText $wsTurf "AB16" "Ronaldo 2`$ captain"
Set-PlainText $wsTurf "AD17" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "Y18" "Benrahma 4`$ captain"
Set-PlainText $wsTurf "AC19" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "AD20" "Ronaldo 2`$ captain"
Set-PlainText $wsTurf "AC21" "Ronaldo 2`$ captain"
Set-PlainText $wsTurf "AC22" "Antonio 4`$ captain"
Set-PlainText $wsTurf "AC23" "Antonio 4`$ captain"
Set-PlainText $wsTurf "AD24" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "AD25" "Jiménez 20`$ captain"
Set-PlainText $wsTurf "AB26" "Ronaldo 2`$ captain"
Set-PlainText $wsTurf "AD27" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "AC28" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "AD29" "Ronaldo 2`$ captain"
Set-PlainText $wsTurf "AD30" "Lukaku 4`$ captain"
Set-PlainText $wsTurf "X31" "Raphinha 6`$ captain"
Set-PlainText $wsTurf "AD32" "Ronaldo 2`$ captain"
Set-PlainText $wsTurf "AC33" "Ronaldo 2`$ captain"

$wsTurf.Range("C39").Select() | Out-Null

# ----- Sheet4 -----
$wsSheet4 = $wb.Worksheets.Item("Sheet4")
Set-PlainText $wsSheet4 "AC2" "Lukaku 4`$ captain"

# ----- RoC sheet -----
$wsRoC = $wb.Worksheets.Item("RoC")

# Tag captain picks for existing rows
Set-PlainText $wsRoC "X2" "Salah 26`$ captain"
Set-PlainText $wsRoC "AC3" "Ronaldo 2`$ captain"
Set-PlainText $wsRoC "AB4" "Lukaku 4`$ captain"
Set-PlainText $wsRoC "AB5" "Ronaldo 2`$ captain"
Set-PlainText $wsRoC "AD6" "Lukaku 4`$ captain"

# New entrant - row 7
$wsRoC.Range("A7").Value = "Camilla Nurkhanov"
$wsRoC.Range("B7").Value = 5306515
Set-PlainText $wsRoC "C7" "40"
Set-PlainText $wsRoC "D7" "410"
Set-PlainText $wsRoC "E7" "1,529,578"
Set-PlainText $wsRoC "G7" "2"
$wsRoC.Range("R7").Value = "Camilla Nurkhanov"
$wsRoC.Range("S7").Value = "Kamilkas United"
$wsRoC.Range("T7").Value = "Martínez 4"
$wsRoC.Range("U7").Value = "Dias 0"
$wsRoC.Range("V7").Value = "Livramento 4"
$wsRoC.Range("W7").Value = "van Dijk 1"
$wsRoC.Range("X7").Value = "Salah 13"
$wsRoC.Range("Y7").Value = "Benrahma 2"
$wsRoC.Range("Z7").Value = "Grealish 2"
$wsRoC.Range("AA7").Value = "Raphinha 3"
$wsRoC.Range("AB7").Value = "Antonio 4`$ captain"
$wsRoC.Range("AC7").Value = "Ings 2"
$wsRoC.Range("AD7").Value = "Jesus 5"
$wsRoC.Range("AE7").Value = "Pickford 3"
$wsRoC.Range("AF7").Value = "Alonso 0"
$wsRoC.Range("AG7").Value = "Gibbs-White 0"
$wsRoC.Range("AH7").Value = "Pollock 0"

# New entrant - row 8
$wsRoC.Range("A8").Value = "Leonardo Lombardi"
$wsRoC.Range("B8").Value = 2249991
Set-PlainText $wsRoC "C8" "46"
Set-PlainText $wsRoC "D8" "400"
Set-PlainText $wsRoC "E8" "1,923,784"
Set-PlainText $wsRoC "G8" "0"
$wsRoC.Range("R8").Value = "Leonardo Lombardi"
$wsRoC.Range("S8").Value = "TheChiellinis"
$wsRoC.Range("T8").Value = "Pickford 3"
$wsRoC.Range("U8").Value = "Shaw 2"
$wsRoC.Range("V8").Value = "Dias 0"
$wsRoC.Range("W8").Value = "Dier 3"
$wsRoC.Range("X8").Value = "Dallas 3"
$wsRoC.Range("Y8").Value = "Jota 2"
$wsRoC.Range("Z8").Value = "Gallagher 2"
$wsRoC.Range("AA8").Value = "Salah 26`$ captain"
$wsRoC.Range("AB8").Value = "Antonio 2"
$wsRoC.Range("AC8").Value = "Ronaldo 1"
$wsRoC.Range("AD8").Value = "Ings 2"
$wsRoC.Range("AE8").Value = "Sánchez 6"
$wsRoC.Range("AF8").Value = "Söyüncü 1"
$wsRoC.Range("AG8").Value = "Allan 2"
$wsRoC.Range("AH8").Value = "Mings 1"

# New entrant - row 9
$wsRoC.Range("A9").Value = "Lorin Minxhozi"
$wsRoC.Range("B9").Value = 2056865
Set-PlainText $wsRoC "C9" "21"
Set-PlainText $wsRoC "D9" "358"
Set-PlainText $wsRoC "E9" "3,717,085"
Set-PlainText $wsRoC "G9" "2"
$wsRoC.Range("R9").Value = "Lorin Minxhozi"
$wsRoC.Range("S9").Value = "Pink Lemonade"
$wsRoC.Range("T9").Value = "Sá 2"
$wsRoC.Range("U9").Value = "Jansson 5"
$wsRoC.Range("V9").Value = "Dias 0"
$wsRoC.Range("W9").Value = "Cash 1"
$wsRoC.Range("X9").Value = "Gallagher 2"
$wsRoC.Range("Y9").Value = "Greenwood 1"
$wsRoC.Range("Z9").Value = "Gray 2"
$wsRoC.Range("AA9").Value = "Torres 0"
$wsRoC.Range("AB9").Value = "Jota 2"
$wsRoC.Range("AC9").Value = "Kane 4`$ captain"
$wsRoC.Range("AD9").Value = "Lukaku 2"
$wsRoC.Range("AE9").Value = "Foster 2"
$wsRoC.Range("AF9").Value = "Cucho 0"
$wsRoC.Range("AG9").Value = "Alexander-Arnold 0"
$wsRoC.Range("AH9").Value = "James 0"

# Stray row 10
$wsRoC.Range("A10").Value = "ddf"

# Make RoC the active sheet / cursor on B9 (mirrors the author reviewing the new row)
$wsRoC.Activate()
$wsRoC.Range("B9").Select() | Out-Null

Write-Output "edit complete"